# Apply updates for "Adding labs 24, 26-29"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update window width on the workbook view
$excel.ActiveWindow.Width = 20730

# Row 3
$ws.Range("B3").Value = 7.38
$ws.Range("C3").Value = 7.14

# Row 4
$ws.Range("B4").Value = 44
$ws.Range("C4").Value = 24

# Row 5
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = -646

# Row 6
$ws.Range("B6").Value = 0.28000000000000003
$ws.Range("C6").Value = 199.68
$ws.Range("D6").Value = 396.99

# Row 7
$ws.Range("B7").Value = 48
$ws.Range("C7").Value = 45

# Row 8
$ws.Range("B8").Value = 1.1000000000000001
$ws.Range("C8").Value = 45.7
$ws.Range("D8").Value = 0

# Row 9
$ws.Range("B9").Value = 1.1000000000000001
$ws.Range("C9").Value = 24
$ws.Range("D9").Value = 0

# Row 10
$ws.Range("D10").Value = 0

# Row 11
$ws.Range("B11").Value = 1.1200000000000001
$ws.Range("C11").Value = 1.64
$ws.Range("D11").Value = 0

# Row 12
$ws.Range("B12").Value = 6.5
$ws.Range("C12").Value = 10.7
$ws.Range("D12").Value = 0
